$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert "COTOVELO 90 GR PVC JS DN 20 (DR 15 MM)" row before current row 5
# (the existing "COTOVELO 90 GR PVC JR DN 25 (DR 20 MM)" row)
$ws.Rows("5:5").Insert()
$ws.Range("A5").Value = "30000732"
$ws.Range("B5").Value = "COTOVELO 90 GR PVC JS DN 20 (DR 15 MM)"

# Step 2: insert "PORCA PVC BRANCA DN 20" row before current row 2
# (the existing "REGISTRO PVC  MARROM  DN 20 MM" row)
$ws.Rows("2:2").Insert()
$ws.Range("A2").Value = "50000333"
$ws.Range("B2").Value = "PORCA PVC BRANCA DN 20"

# Step 3: insert "TE PVC JS DN 20 X 20" row before current row 16
# (the existing "TE PVC JS DN 25X25" row)
$ws.Rows("16:16").Insert()
$ws.Range("A16").Value = "30003758"
$ws.Range("B16").Value = "TE PVC JS DN 20 X 20"

# Step 4: insert "TUBO PVC RIGIDO PB JS DN 25 (DR 20 MM)" row before current row 20
# (the existing "TUBO PVC RIG PP JR DN 25 (DR 20 MM)" row)
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = "30000880"
$ws.Range("B20").Value = "TUBO PVC RIGIDO PB JS DN 25 (DR 20 MM)"

# Update the sheet selection to match the final cursor position
$ws.Range("B28").Select()
